# Highlight the title on slide 1: "Modeling, Simulation, and Analysis (MSA)"
# Everything up to and including "MSA" gets a yellow highlight; the closing
# parenthesis ")" is left un-highlighted (matching the target OOXML, which
# splits the single run into two runs).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the title shape (the "ctrTitle" placeholder) on the slide by matching
# its current text instead of a hard-coded shape index, so the script is
# resilient to shape ordering.
$targetText = "Modeling, Simulation, and Analysis (MSA)"
$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -eq $targetText) {
        $titleShape = $candidate
        break
    }
}
if ($titleShape -eq $null) {
    $titleShape = $s.Shapes.Item(2)
}

$tr = $titleShape.TextFrame.TextRange
$fullText = $tr.Text
$splitPoint = $fullText.Length - 1   # length of "Modeling, Simulation, and Analysis (MSA" (everything but the final ")")

# Sub-range covering "Modeling, Simulation, and Analysis (MSA" -> gets the highlight.
$highlightRange = $tr.Characters(1, $splitPoint)
$highlightRange.Font.Highlight.RGB = 65535   # RGB(255,255,0) = yellow, stored as &H00FFFF00 -> 65535

